$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.299608333333333
$ws.Range("H2").Value = 6.898825
$ws.Range("I2").Value = 0.004282915438164647
$ws.Range("J2").Value = 0.004282915438164647
$ws.Range("M2").Value = 7.085731
$ws.Range("N2").Value = 21.257193
$ws.Range("O2").Value = 0.1632021414334214
$ws.Range("P2").Value = 0.1632021414334214
$ws.Range("Q2").Value = 16.29440605535833
$ws.Range("R2").Value = 146.649654498225
$ws.Range("S2").Value = 0.0006989809710867308
$ws.Range("T2").Value = 0.0006989809710867308

$ws.Range("G3").Value = 2.299608333333333
$ws.Range("H3").Value = 6.898825
$ws.Range("I3").Value = 0.004282915438164647
$ws.Range("J3").Value = 0.004282915438164647
$ws.Range("O3").Value = 0.2146673930709577
$ws.Range("P3").Value = 0.2146673930709577
$ws.Range("Q3").Value = 21.43279272453889
$ws.Range("R3").Value = 192.89513452085
$ws.Range("S3").Value = 0.0009194022918541633
$ws.Range("T3").Value = 0.0009194022918541634

$ws.Range("G4").Value = 2.299608333333333
$ws.Range("H4").Value = 6.898825
$ws.Range("I4").Value = 0.004282915438164647
$ws.Range("J4").Value = 0.004282915438164647
$ws.Range("M4").Value = 14.22029466666667
$ws.Range("N4").Value = 42.660884
$ws.Range("O4").Value = 0.3275290215525062
$ws.Range("P4").Value = 0.3275290215525062
$ws.Range("Q4").Value = 32.70110811792222
$ws.Range("R4").Value = 294.3099730613
$ws.Range("S4").Value = 0.00140277910285419
$ws.Range("T4").Value = 0.00140277910285419

$ws.Range("G5").Value = 2.299608333333333
$ws.Range("H5").Value = 6.898825
$ws.Range("I5").Value = 0.004282915438164647
$ws.Range("J5").Value = 0.004282915438164647
$ws.Range("M5").Value = 1.993361666666667
$ws.Range("N5").Value = 5.980085
$ws.Range("O5").Value = 0.04591211445245296
$ws.Range("P5").Value = 0.04591211445245296
$ws.Range("Q5").Value = 4.58395110001389
$ws.Range("R5").Value = 41.255559900125
$ws.Range("S5").Value = 0.000196637703787193
$ws.Range("T5").Value = 0.000196637703787193

$ws.Range("G6").Value = 2.299608333333333
$ws.Range("H6").Value = 6.898825
$ws.Range("I6").Value = 0.004282915438164647
$ws.Range("J6").Value = 0.004282915438164647
$ws.Range("M6").Value = 10.79731966666667
$ws.Range("N6").Value = 32.391959
$ws.Range("O6").Value = 0.2486893294906617
$ws.Range("P6").Value = 0.2486893294906617
$ws.Range("Q6").Value = 24.82960628313056
$ws.Range("R6").Value = 223.466456548175
$ws.Range("S6").Value = 0.00106511536858237
$ws.Range("T6").Value = 0.00106511536858237

$ws.Range("I7").Value = 0.01154106748213608
$ws.Range("J7").Value = 0.01154106748213608
$ws.Range("M7").Value = 7.085731
$ws.Range("N7").Value = 21.257193
$ws.Range("O7").Value = 0.1632021414334214
$ws.Range("P7").Value = 0.1632021414334214
$ws.Range("Q7").Value = 43.90813747814834
$ws.Range("R7").Value = 395.173237303335
$ws.Range("S7").Value = 0.001883526927512233
$ws.Range("T7").Value = 0.001883526927512233

$ws.Range("I8").Value = 0.01154106748213608
$ws.Range("J8").Value = 0.01154106748213608
$ws.Range("O8").Value = 0.2146673930709577
$ws.Range("P8").Value = 0.2146673930709577
$ws.Range("S8").Value = 0.002477490869646153
$ws.Range("T8").Value = 0.002477490869646153

$ws.Range("I9").Value = 0.01154106748213608
$ws.Range("J9").Value = 0.01154106748213608
$ws.Range("M9").Value = 14.22029466666667
$ws.Range("N9").Value = 42.660884
$ws.Range("O9").Value = 0.3275290215525062
$ws.Range("P9").Value = 0.3275290215525062
$ws.Range("Q9").Value = 88.11887626044222
$ws.Range("R9").Value = 793.06988634398
$ws.Range("S9").Value = 0.003780034540095475
$ws.Range("T9").Value = 0.003780034540095475

$ws.Range("I10").Value = 0.01154106748213608
$ws.Range("J10").Value = 0.01154106748213608
$ws.Range("M10").Value = 1.993361666666667
$ws.Range("N10").Value = 5.980085
$ws.Range("O10").Value = 0.04591211445245296
$ws.Range("P10").Value = 0.04591211445245296
$ws.Range("Q10").Value = 12.35226091756389
$ws.Range("R10").Value = 111.170348258075
$ws.Range("S10").Value = 0.0005298748111433147
$ws.Range("T10").Value = 0.0005298748111433147

$ws.Range("I11").Value = 0.01154106748213608
$ws.Range("J11").Value = 0.01154106748213608
$ws.Range("M11").Value = 10.79731966666667
$ws.Range("N11").Value = 32.391959
$ws.Range("O11").Value = 0.2486893294906617
$ws.Range("P11").Value = 0.2486893294906617
$ws.Range("Q11").Value = 66.90773278290057
$ws.Range("R11").Value = 602.1695950461051
$ws.Range("S11").Value = 0.0028701403337389
$ws.Range("T11").Value = 0.0028701403337389

$ws.Range("G12").Value = 265.3156126666667
$ws.Range("H12").Value = 795.9468380000001
$ws.Range("I12").Value = 0.4941382047563949
$ws.Range("J12").Value = 0.4941382047563948
$ws.Range("M12").Value = 7.085731
$ws.Range("N12").Value = 21.257193
$ws.Range("O12").Value = 0.1632021414334214
$ws.Range("P12").Value = 0.1632021414334214
$ws.Range("Q12").Value = 1879.955061456193
$ws.Range("R12").Value = 16919.59555310574
$ws.Range("S12").Value = 0.0806444131803101
$ws.Range("T12").Value = 0.08064441318031008

$ws.Range("G13").Value = 265.3156126666667
$ws.Range("H13").Value = 795.9468380000001
$ws.Range("I13").Value = 0.4941382047563949
$ws.Range("J13").Value = 0.4941382047563948
$ws.Range("O13").Value = 0.2146673930709577
$ws.Range("P13").Value = 0.2146673930709577
$ws.Range("Q13").Value = 2472.792627528041
$ws.Range("R13").Value = 22255.13364775237
$ws.Range("S13").Value = 0.1060753602318184
$ws.Range("T13").Value = 0.1060753602318184

$ws.Range("G14").Value = 265.3156126666667
$ws.Range("H14").Value = 795.9468380000001
$ws.Range("I14").Value = 0.4941382047563949
$ws.Range("J14").Value = 0.4941382047563948
$ws.Range("M14").Value = 14.22029466666667
$ws.Range("N14").Value = 42.660884
$ws.Range("O14").Value = 0.3275290215525062
$ws.Range("P14").Value = 0.3275290215525062
$ws.Range("Q14").Value = 3772.8661917872
$ws.Range("R14").Value = 33955.79572608479
$ws.Range("S14").Value = 0.161844602715574
$ws.Range("T14").Value = 0.1618446027155739

$ws.Range("G15").Value = 265.3156126666667
$ws.Range("H15").Value = 795.9468380000001
$ws.Range("I15").Value = 0.4941382047563949
$ws.Range("J15").Value = 0.4941382047563948
$ws.Range("M15").Value = 1.993361666666667
$ws.Range("N15").Value = 5.980085
$ws.Range("O15").Value = 0.04591211445245296
$ws.Range("P15").Value = 0.04591211445245296
$ws.Range("Q15").Value = 528.8699718579145
$ws.Range("R15").Value = 4759.829746721231
$ws.Range("S15").Value = 0.02268692981210524
$ws.Range("T15").Value = 0.02268692981210524

$ws.Range("G16").Value = 265.3156126666667
$ws.Range("H16").Value = 795.9468380000001
$ws.Range("I16").Value = 0.4941382047563949
$ws.Range("J16").Value = 0.4941382047563948
$ws.Range("M16").Value = 10.79731966666667
$ws.Range("N16").Value = 32.391959
$ws.Range("O16").Value = 0.2486893294906617
$ws.Range("P16").Value = 0.2486893294906617
$ws.Range("Q16").Value = 2864.697482519516
$ws.Range("R16").Value = 25782.27734267564
$ws.Range("S16").Value = 0.1228868988165871
$ws.Range("T16").Value = 0.1228868988165871

$ws.Range("G17").Value = 16.35749033333333
$ws.Range("H17").Value = 49.072471
$ws.Range("I17").Value = 0.03046507827561751
$ws.Range("J17").Value = 0.0304650782756175
$ws.Range("M17").Value = 7.085731
$ws.Range("N17").Value = 21.257193
$ws.Range("O17").Value = 0.1632021414334214
$ws.Range("P17").Value = 0.1632021414334214
$ws.Range("Q17").Value = 115.9047763371003
$ws.Range("R17").Value = 1043.142987033903
$ws.Range("S17").Value = 0.004971966013517582
$ws.Range("T17").Value = 0.004971966013517582

$ws.Range("G18").Value = 16.35749033333333
$ws.Range("H18").Value = 49.072471
$ws.Range("I18").Value = 0.03046507827561751
$ws.Range("J18").Value = 0.0304650782756175
$ws.Range("O18").Value = 0.2146673930709577
$ws.Range("P18").Value = 0.2146673930709577
$ws.Range("Q18").Value = 152.4549614498042
$ws.Range("R18").Value = 1372.094653048238
$ws.Range("S18").Value = 0.006539858933129477
$ws.Range("T18").Value = 0.006539858933129477

$ws.Range("G19").Value = 16.35749033333333
$ws.Range("H19").Value = 49.072471
$ws.Range("I19").Value = 0.03046507827561751
$ws.Range("J19").Value = 0.0304650782756175
$ws.Range("M19").Value = 14.22029466666667
$ws.Range("N19").Value = 42.660884
$ws.Range("O19").Value = 0.3275290215525062
$ws.Range("P19").Value = 0.3275290215525062
$ws.Range("Q19").Value = 232.6083325471515
$ws.Range("R19").Value = 2093.474992924364
$ws.Range("S19").Value = 0.009978197279133514
$ws.Range("T19").Value = 0.009978197279133513

$ws.Range("G20").Value = 16.35749033333333
$ws.Range("H20").Value = 49.072471
$ws.Range("I20").Value = 0.03046507827561751
$ws.Range("J20").Value = 0.0304650782756175
$ws.Range("M20").Value = 1.993361666666667
$ws.Range("N20").Value = 5.980085
$ws.Range("O20").Value = 0.04591211445245296
$ws.Range("P20").Value = 0.04591211445245296
$ws.Range("Q20").Value = 32.60639419333722
$ws.Range("R20").Value = 293.457547740035
$ws.Range("S20").Value = 0.001398716160593089
$ws.Range("T20").Value = 0.001398716160593089

$ws.Range("G21").Value = 16.35749033333333
$ws.Range("H21").Value = 49.072471
$ws.Range("I21").Value = 0.03046507827561751
$ws.Range("J21").Value = 0.0304650782756175
$ws.Range("M21").Value = 10.79731966666667
$ws.Range("N21").Value = 32.391959
$ws.Range("O21").Value = 0.2486893294906617
$ws.Range("P21").Value = 0.2486893294906617
$ws.Range("Q21").Value = 176.6170520734099
$ws.Range("R21").Value = 1589.553468660689
$ws.Range("S21").Value = 0.007576339889243842
$ws.Range("T21").Value = 0.007576339889243841

$ws.Range("G22").Value = 246.7565153333333
$ws.Range("H22").Value = 740.269546
$ws.Range("I22").Value = 0.459572734047687
$ws.Range("J22").Value = 0.459572734047687
$ws.Range("M22").Value = 7.085731
$ws.Range("N22").Value = 21.257193
$ws.Range("O22").Value = 0.1632021414334214
$ws.Range("P22").Value = 0.1632021414334214
$ws.Range("Q22").Value = 1748.450290149375
$ws.Range("R22").Value = 15736.05261134438
$ws.Range("S22").Value = 0.07500325434099478
$ws.Range("T22").Value = 0.07500325434099478

$ws.Range("G23").Value = 246.7565153333333
$ws.Range("H23").Value = 740.269546
$ws.Range("I23").Value = 0.459572734047687
$ws.Range("J23").Value = 0.459572734047687
$ws.Range("O23").Value = 0.2146673930709577
$ws.Range("P23").Value = 0.2146673930709577
$ws.Range("Q23").Value = 2299.818264661954
$ws.Range("R23").Value = 20698.36438195759
$ws.Range("S23").Value = 0.09865528074450954
$ws.Range("T23").Value = 0.09865528074450954

$ws.Range("G24").Value = 246.7565153333333
$ws.Range("H24").Value = 740.269546
$ws.Range("I24").Value = 0.459572734047687
$ws.Range("J24").Value = 0.459572734047687
$ws.Range("M24").Value = 14.22029466666667
$ws.Range("N24").Value = 42.660884
$ws.Range("O24").Value = 0.3275290215525062
$ws.Range("P24").Value = 0.3275290215525062
$ws.Range("Q24").Value = 3508.950358959852
$ws.Range("R24").Value = 31580.55323063866
$ws.Range("S24").Value = 0.1505234079148491
$ws.Range("T24").Value = 0.1505234079148491

$ws.Range("G25").Value = 246.7565153333333
$ws.Range("H25").Value = 740.269546
$ws.Range("I25").Value = 0.459572734047687
$ws.Range("J25").Value = 0.459572734047687
$ws.Range("M25").Value = 1.993361666666667
$ws.Range("N25").Value = 5.980085
$ws.Range("O25").Value = 0.04591211445245296
$ws.Range("P25").Value = 0.04591211445245296
$ws.Range("Q25").Value = 491.8749786657123
$ws.Range("R25").Value = 4426.87480799141
$ws.Range("S25").Value = 0.02109995596482413
$ws.Range("T25").Value = 0.02109995596482413

$ws.Range("G26").Value = 246.7565153333333
$ws.Range("H26").Value = 740.269546
$ws.Range("I26").Value = 0.459572734047687
$ws.Range("J26").Value = 0.459572734047687
$ws.Range("M26").Value = 10.79731966666667
$ws.Range("N26").Value = 32.391959
$ws.Range("O26").Value = 0.2486893294906617
$ws.Range("P26").Value = 0.2486893294906617
$ws.Range("Q26").Value = 2664.308975886735
$ws.Range("R26").Value = 23978.78078298062
$ws.Range("S26").Value = 0.1142908350825095
$ws.Range("T26").Value = 0.1142908350825095
